$d = $word.ActiveDocument

# Delete the "Ver no Jupiter Salvar em pdf Salvar em docx" paragraph and the
# following "(c) 2020 ... Creative Commons Attribution" paragraph (the site
# footer that was dropped from the source page), leaving everything else -
# including the single blank paragraph before them - untouched.

$jupiterText = "Ver no Jupiter Salvar em pdf Salvar em docx"
$copyrightText = [char]0xA9 + " 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution"

for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($t -eq $jupiterText -or $t -eq $copyrightText) {
        $p.Range.Delete()
    }
}
